# --- Part 1: Add new "2022-Q1" worksheet, positioned after "2021-Q4" and before "总计" ---
$wb = $excel.ActiveWorkbook

$srcQ4 = $wb.Worksheets.Item("2021-Q4")
$srcQ4.Copy($null, $srcQ4)
$newSheet = $wb.Worksheets.Item($srcQ4.Index + 1)
$newSheet.Name = "2022-Q1"

# Extend column-A / row formatting (style copied from the source sheet's last
# data row) down through row 15 so the 11 additional data rows inherit the
# same look (bordered, bold, centered index column) as the existing rows.
$newSheet.Range("A4:H4").Copy()
$newSheet.Range("A5:H15").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

# --- Part 2: populate the 14 data rows (index, code, name, scale, position,
# ratio, market value, rank) ---
$ws = $newSheet
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'000404"
$ws.Cells.Item(2,3).Value = "'易方达新兴成长灵活配置混合"
$ws.Cells.Item(2,4).Value = "'51.67"
$ws.Cells.Item(2,5).Value = "'91.13"
$ws.Cells.Item(2,6).Value = "'4.93"
$ws.Cells.Item(2,7).Value = "'2.5473"
$ws.Cells.Item(2,8).Value = 9
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'009808"
$ws.Cells.Item(3,3).Value = "'易方达创新成长混合"
$ws.Cells.Item(3,4).Value = "'33.50"
$ws.Cells.Item(3,5).Value = "'93.21"
$ws.Cells.Item(3,6).Value = "'4.55"
$ws.Cells.Item(3,7).Value = "'1.5242"
$ws.Cells.Item(3,8).Value = 10
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'213001"
$ws.Cells.Item(4,3).Value = "'宝盈鸿利收益灵活配置混合A"
$ws.Cells.Item(4,4).Value = "'17.98"
$ws.Cells.Item(4,5).Value = "'90.37"
$ws.Cells.Item(4,6).Value = "'3.75"
$ws.Cells.Item(4,7).Value = "'0.6742"
$ws.Cells.Item(4,8).Value = 10
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'012850"
$ws.Cells.Item(5,3).Value = "'中融低碳经济3个月持有期混合型证券投资基金A"
$ws.Cells.Item(5,4).Value = "'8.67"
$ws.Cells.Item(5,5).Value = "'65.64"
$ws.Cells.Item(5,6).Value = "'3.00"
$ws.Cells.Item(5,7).Value = "'0.2601"
$ws.Cells.Item(5,8).Value = 9
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'010751"
$ws.Cells.Item(6,3).Value = "'宝盈优质成长混合A"
$ws.Cells.Item(6,4).Value = "'5.64"
$ws.Cells.Item(6,5).Value = "'92.80"
$ws.Cells.Item(6,6).Value = "'4.33"
$ws.Cells.Item(6,7).Value = "'0.2442"
$ws.Cells.Item(6,8).Value = 8
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'001543"
$ws.Cells.Item(7,3).Value = "'宝盈新锐灵活配置混合A"
$ws.Cells.Item(7,4).Value = "'3.21"
$ws.Cells.Item(7,5).Value = "'93.26"
$ws.Cells.Item(7,6).Value = "'4.87"
$ws.Cells.Item(7,7).Value = "'0.1563"
$ws.Cells.Item(7,8).Value = 7
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'013956"
$ws.Cells.Item(8,3).Value = "'华商医药消费精选混合A"
$ws.Cells.Item(8,4).Value = "'5.01"
$ws.Cells.Item(8,5).Value = "'25.45"
$ws.Cells.Item(8,6).Value = "'1.82"
$ws.Cells.Item(8,7).Value = "'0.0912"
$ws.Cells.Item(8,8).Value = 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'012851"
$ws.Cells.Item(9,3).Value = "'中融低碳经济3个月持有期混合型证券投资基金C"
$ws.Cells.Item(9,4).Value = "'1.49"
$ws.Cells.Item(9,5).Value = "'65.64"
$ws.Cells.Item(9,6).Value = "'3.00"
$ws.Cells.Item(9,7).Value = "'0.0447"
$ws.Cells.Item(9,8).Value = 9
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'010752"
$ws.Cells.Item(10,3).Value = "'宝盈优质成长混合C"
$ws.Cells.Item(10,4).Value = "'0.78"
$ws.Cells.Item(10,5).Value = "'92.80"
$ws.Cells.Item(10,6).Value = "'4.33"
$ws.Cells.Item(10,7).Value = "'0.0338"
$ws.Cells.Item(10,8).Value = 8
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'005571"
$ws.Cells.Item(11,3).Value = "'中银证券新能源灵活配置混合A"
$ws.Cells.Item(11,4).Value = "'0.91"
$ws.Cells.Item(11,5).Value = "'90.25"
$ws.Cells.Item(11,6).Value = "'3.43"
$ws.Cells.Item(11,7).Value = "'0.0312"
$ws.Cells.Item(11,8).Value = 9
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'007581"
$ws.Cells.Item(12,3).Value = "'宝盈鸿利收益灵活配置混合C"
$ws.Cells.Item(12,4).Value = "'0.73"
$ws.Cells.Item(12,5).Value = "'90.37"
$ws.Cells.Item(12,6).Value = "'3.75"
$ws.Cells.Item(12,7).Value = "'0.0274"
$ws.Cells.Item(12,8).Value = 10
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'007578"
$ws.Cells.Item(13,3).Value = "'宝盈新锐灵活配置混合C"
$ws.Cells.Item(13,4).Value = "'0.20"
$ws.Cells.Item(13,5).Value = "'93.26"
$ws.Cells.Item(13,6).Value = "'4.87"
$ws.Cells.Item(13,7).Value = "'0.0097"
$ws.Cells.Item(13,8).Value = 7
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'005572"
$ws.Cells.Item(14,3).Value = "'中银证券新能源灵活配置混合C"
$ws.Cells.Item(14,4).Value = "'0.28"
$ws.Cells.Item(14,5).Value = "'90.25"
$ws.Cells.Item(14,6).Value = "'3.43"
$ws.Cells.Item(14,7).Value = "'0.0096"
$ws.Cells.Item(14,8).Value = 9
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "'013957"
$ws.Cells.Item(15,3).Value = "'华商医药消费精选混合C"
$ws.Cells.Item(15,4).Value = "'0.39"
$ws.Cells.Item(15,5).Value = "'25.45"
$ws.Cells.Item(15,6).Value = "'1.82"
$ws.Cells.Item(15,7).Value = "'0.0071"
$ws.Cells.Item(15,8).Value = 9


# Columns B:G hold text-typed figures (fund codes / names / percentages as
# strings, matching the source data) -- the leading "'" above forces text
# entry so things like "000404" keep their leading zeros instead of being
# parsed as numbers. ClearFormats() strips the resulting quote-prefix style
# bit so the cells fall back to the workbook's default (unstyled) format,
# matching the sibling quarter sheets.
$ws.Range("B2:G15").ClearFormats()

[void]$newSheet.Range("A1").Select()

# --- Part 3: update the "总计" (total) summary sheet with the new quarter ---
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

# Pull the A-column index style (bordered/bold/centered) from the row below
# onto the freshly inserted row, then clear the stray format Excel gave the
# rest of the new row during the insert.
$totalWs.Cells.Item(3,1).Copy()
$totalWs.Cells.Item(2,1).PasteSpecial(-4122)
$totalWs.Application.CutCopyMode = $false
$totalWs.Range("B2:D2").ClearFormats()

$totalWs.Cells.Item(2,1).Value = 0
$totalWs.Cells.Item(2,2).Value = "2022-Q1"
$totalWs.Cells.Item(2,3).Value = 14
$totalWs.Cells.Item(2,4).Value = 5.66

# Renumber the index column for the rows that shifted down
$totalWs.Cells.Item(3,1).Value = 1
$totalWs.Cells.Item(4,1).Value = 2
$totalWs.Cells.Item(5,1).Value = 3

[void]$totalWs.Range("A1").Select()

# Restore the originally active tab (first sheet) so the only visible change
# is the new quarter's data, not which tab happens to be selected.
[void]$wb.Worksheets.Item(1).Activate()
